$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (date) and Volumen values between row 2 and row 4
$ws.Range("D2").Value = 44874
$ws.Range("M2").Value = 67

$ws.Range("D4").Value = 44875
$ws.Range("M4").Value = 50
